$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel stores them as text (matching the source data which is all text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '66.267.40'
$ws.Range('E2').Value = '  +5.96%  '
$ws.Range('D3').Value = '3.004.92'
$ws.Range('E3').Value = '  +3.18%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '582.54'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('D6').Value = '163.64'
$ws.Range('E6').Value = '  +12.91%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.518'
$ws.Range('E8').Value = '  +3.37%  '
$ws.Range('D9').Value = '3.000.61'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').Value = '6.60'
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('E12').Value = '  +5.35%  '
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  +6.09%  '
$ws.Range('D14').Value = '34.76'
$ws.Range('E14').Value = '  +5.63%  '
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '66.239.17'
$ws.Range('E16').Value = '  +5.99%  '
$ws.Range('D17').Value = '3.502.99'
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('D18').Value = '6.93'
$ws.Range('E18').Value = '  +4.43%  '
$ws.Range('D19').Value = '3.005.93'
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('D20').Value = '455.52'
$ws.Range('E20').Value = '  +5.49%  '
$ws.Range('D21').Value = '13.86'
$ws.Range('E21').Value = '  +5.29%  '
$ws.Range('E22').Value = '  +3.94%  '
$ws.Range('D23').Value = '7.36'
$ws.Range('E23').Value = '  +6.65%  '
$ws.Range('D24').Value = '82.39'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').Value = '2.30'
$ws.Range('E25').Value = '  +13.53%  '
$ws.Range('E26').Value = '  +2.85%  '
$ws.Range('D27').Value = '10.48'
$ws.Range('E27').Value = '  +4.69%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '8.16'
$ws.Range('E29').Value = '  +15.47%  '
$ws.Range('D30').Value = '2.39'
$ws.Range('E30').Value = '  +17.72%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.64'
$ws.Range('E31').Value = '  +5.75%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0000105'
$ws.Range('E32').Value = '  -5.06%  '
$ws.Range('D33').Value = '27.30'
$ws.Range('E33').Value = '  +5.73%  '
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  +3.87%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').Value = '0.992'
$ws.Range('E36').Value = '  +3.97%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = '5.85'
$ws.Range('E37').Value = '  +8.03%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '49.87'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '2.05'
$ws.Range('E39').Value = '  +7.59%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').Value = '0.310'
$ws.Range('E41').Value = '  +15.63%  '
$ws.Range('D42').Value = '0.122'
$ws.Range('E42').Value = '  +7.05%  '
$ws.Range('D43').Value = '43.71'
$ws.Range('E43').Value = '  +6.02%  '
$ws.Range('E44').Value = '  +4.27%  '
$ws.Range('D45').Value = '403.63'
$ws.Range('E45').Value = '  +14.64%  '
$ws.Range('D46').Value = '0.0359'
$ws.Range('E46').Value = '  +5.92%  '
$ws.Range('D47').Value = '2.792.41'
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('D48').Value = '134.34'
$ws.Range('E48').Value = '  +0.71%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').Value = '23.89'
$ws.Range('E50').Value = '  +11.29%  '
$ws.Range('E51').Value = '  +4.12%  '
